$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.554.22"
$ws.Range("E2").Value = "  +4.20%  "
$ws.Range("D3").Value = "3.252.41"
$ws.Range("E3").Value = "  +2.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.14"
$ws.Range("E5").Value = "  +1.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.40"
$ws.Range("E6").Value = "  +7.13%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -3.64%  "
$ws.Range("D9").Value = "3.251.45"
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("E10").Value = "  +6.16%  "
$ws.Range("E11").Value = "  +3.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.415"
$ws.Range("E12").Value = "  +5.08%  "
$ws.Range("D13").Value = "3.818.19"
$ws.Range("E13").Value = "  +2.80%  "
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.60"
$ws.Range("E15").Value = "  +5.40%  "
$ws.Range("D16").Value = "67.495.09"
$ws.Range("E16").Value = "  +4.18%  "
$ws.Range("E17").Value = "  +3.37%  "
$ws.Range("D18").Value = "3.245.87"
$ws.Range("E18").Value = "  +2.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.85"
$ws.Range("E19").Value = "  +2.00%  "
$ws.Range("E20").Value = "  +5.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.25"
$ws.Range("E21").Value = "  +5.65%  "
$ws.Range("E22").Value = "  +4.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.39"
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.512"
$ws.Range("E25").Value = "  +2.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.61"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("E28").Value = "  +2.72%  "
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("E30").Value = "  +8.43%  "
$ws.Range("E31").Value = "  +3.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.71"
$ws.Range("E32").Value = "  +3.25%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.28"
$ws.Range("E34").Value = "  +5.93%  "
$ws.Range("E35").Value = "  +4.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.36"
$ws.Range("E36").Value = "  +5.48%  "
$ws.Range("E37").Value = "  +4.09%  "
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("E39").Value = "  +5.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.85"
$ws.Range("E40").Value = "  +13.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.87"
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("E42").Value = "  +10.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.60"
$ws.Range("E43").Value = "  +4.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "358.35"
$ws.Range("E44").Value = "  +10.69%  "
$ws.Range("D45").Value = "2.727.70"
$ws.Range("E45").Value = "  +2.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.45"
$ws.Range("E46").Value = "  +5.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.89"
$ws.Range("E47").Value = "  +3.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0680"
$ws.Range("E48").Value = "  +3.21%  "
$ws.Range("E49").Value = "  +2.31%  "
$ws.Range("E50").Value = "  +6.49%  "
$ws.Range("E51").Value = "  -0.63%  "
